# PrecioFrutaHortalizas - Pepino ensalada: add a new weekly price record.
# A brand-new observation is inserted as row 614 (pushing the existing
# rows 614-681 down to 615-682), reflecting the "Fruta / hortaliza, semanal"
# update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a fresh row at position 614.
$ws.Rows("614:614").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(614, 1).Value2  = 6
$ws.Cells.Item(614, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(614, 3).Value2  = "Metropolitana"
$ws.Cells.Item(614, 4).Value2  = 45194
$ws.Cells.Item(614, 5).Value2  = 13
$ws.Cells.Item(614, 6).Value2  = 100112043
$ws.Cells.Item(614, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(614, 8).Value2  = "Sin especificar"
$ws.Cells.Item(614, 9).Value2  = "Primera"
$ws.Cells.Item(614, 10).Value2 = 800
$ws.Cells.Item(614, 11).Value2 = 11000
$ws.Cells.Item(614, 12).Value2 = 12000
$ws.Cells.Item(614, 13).Value2 = 11562
$ws.Cells.Item(614, 14).Value2 = "$/caja 70 unidades"
$ws.Cells.Item(614, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(614, 16).Value2 = 165
$ws.Cells.Item(614, 17).Value2 = 70
$ws.Cells.Item(614, 18).Value2 = "Hortaliza"
